$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.228.25"
$ws.Range("E2").Value = "  -6.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.55"
$ws.Range("E3").Value = "  -4.07%  "

$ws.Range("E4").Value = "  +0.57%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.83"
$ws.Range("E5").Value = "  -3.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5067"
$ws.Range("E6").Value = "  -12.49%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2657"
$ws.Range("E8").Value = "  -2.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06342"
$ws.Range("E9").Value = "  -4.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.55"
$ws.Range("E10").Value = "  -7.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07375"
$ws.Range("E11").Value = "  -2.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.677.59"
$ws.Range("E12").Value = "  -3.57%  "

$ws.Range("E13").Value = "  -3.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5800"
$ws.Range("E14").Value = "  -3.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.894.62"
$ws.Range("E15").Value = "  -4.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008521"
$ws.Range("E16").Value = "  -2.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.80"
$ws.Range("E17").Value = "  -13.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.172.45"
$ws.Range("E18").Value = "  -6.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.931"
$ws.Range("E19").Value = "  -7.29%  "

$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.82"
$ws.Range("E21").Value = "  -4.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "189.02"
$ws.Range("E22").Value = "  -8.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.195"

$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.84"
$ws.Range("E25").Value = "  -4.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.694"
$ws.Range("E26").Value = "  -4.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1170"
$ws.Range("E27").Value = "  -4.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.69"
$ws.Range("E28").Value = "  -2.94%  "

$ws.Range("E29").Value = "  -5.97%  "

$ws.Range("E30").Value = "  -8.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.322"
$ws.Range("E31").Value = "  -5.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.530"
$ws.Range("E32").Value = "  -5.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.514"
$ws.Range("E33").Value = "  -6.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.632"
$ws.Range("E34").Value = "  -2.37%  "

$ws.Range("E35").Value = "  -2.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5979"
$ws.Range("E36").Value = "  -6.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.359"
$ws.Range("E37").Value = "  -2.45%  "

$ws.Range("E39").Value = "  -3.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.003"
$ws.Range("E40").Value = "  -2.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.074.13"
$ws.Range("E41").Value = "  -4.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8603"
$ws.Range("E42").Value = "  -1.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.45"
$ws.Range("E44").Value = "  -0.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.817.83"
$ws.Range("E45").Value = "  -3.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000111"
$ws.Range("E46").Value = "  +2.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.67"
$ws.Range("E47").Value = "  -6.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  +0.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.058"
$ws.Range("E49").Value = "  -2.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4301"
$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05181"
